$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames
$ws.Range("A1").Value = "Abbr."
$ws.Range("B1").Value = "Subjects"

# Project "subject" column renames (merged cells, value lives in the top-left cell)
$ws.Range("B2").Value = "Apache Commons Lang"
$ws.Range("B9").Value = "UAA"
$ws.Range("B13").Value = "Sql Parser"
$ws.Range("B18").Value = "Joda-Time"
$ws.Range("B22").Value = "Message Pack for Java"
$ws.Range("B29").Value = "Java APNS"
$ws.Range("B35").Value = "Linear Algebra for Java"
$ws.Range("B42").Value = "Wire Mobile Protocol Buffers"

# Update the selected range shown in the saved view
$ws.Range("B2:B8").Select()
